# Insert two new data rows at row 120 (pushing the existing rows 120-198
# down to 122-200) and populate them with the new record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 120; existing rows shift down.
$ws.Rows("120:121").Insert()

# New row 120 (Primera)
$ws.Cells.Item(120, 1).Value = 9
$ws.Cells.Item(120, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(120, 3).Value = "Metropolitana"
$ws.Cells.Item(120, 4).Value = 44518
$ws.Cells.Item(120, 5).Value = 13
$ws.Cells.Item(120, 6).Value = 100112001
$ws.Cells.Item(120, 7).Value = "Berenjena"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 106
$ws.Cells.Item(120, 11).Value = 6000
$ws.Cells.Item(120, 12).Value = 7000
$ws.Cells.Item(120, 13).Value = 6500
$ws.Cells.Item(120, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(120, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(120, 16).Value = 130
$ws.Cells.Item(120, 17).Value = 50
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# New row 121 (Segunda)
$ws.Cells.Item(121, 1).Value = 9
$ws.Cells.Item(121, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(121, 3).Value = "Metropolitana"
$ws.Cells.Item(121, 4).Value = 44518
$ws.Cells.Item(121, 5).Value = 13
$ws.Cells.Item(121, 6).Value = 100112001
$ws.Cells.Item(121, 7).Value = "Berenjena"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Segunda"
$ws.Cells.Item(121, 10).Value = 52
$ws.Cells.Item(121, 11).Value = 5000
$ws.Cells.Item(121, 12).Value = 5000
$ws.Cells.Item(121, 13).Value = 5000
$ws.Cells.Item(121, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 16).Value = 50
$ws.Cells.Item(121, 17).Value = 100
$ws.Cells.Item(121, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date/time number format used by
# the rest of column D (style index carried over automatically from the
# inserted rows, but set explicitly as a safety net).
$ws.Range("D120:D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
